$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D width: 30.85546875 -> 104 (ColumnWidth input compensates for the
# engine's +5/6 conversion offset so the stored xlsx width lands on exactly 104) ---
$ws.Range("D1").ColumnWidth = 619/6

# --- Rows that gain a "Marked" flag in column F (reusing the existing shared
# string "Marked"), and pick up the same formatting (fill + left align) already
# used on similarly-flagged rows (e.g. row 2/3/7/33). We copy formats from a
# donor row so the engine reuses the existing style indices instead of minting
# new ones. ---

# Row 4
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C4:E4").PasteSpecial(-4122)
$ws.Range("F4").Value = "Marked"

# Row 6
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C6:E6").PasteSpecial(-4122)
$ws.Range("F6").Value = "Marked"

# Row 28
$ws.Range("B2").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C28:E28").PasteSpecial(-4122)
$ws.Range("F28").Value = "Marked"

# Row 29
$ws.Range("B2").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C29:E29").PasteSpecial(-4122)
$ws.Range("F29").Value = "Marked"

# Row 34
$ws.Range("B2").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C34:E34").PasteSpecial(-4122)
$ws.Range("F34").Value = "Marked"

$excel.CutCopyMode = $false

# --- View state: scroll down a bit and leave the cursor just past the last
# data row. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("D42").Select()
